# "Coeficiente PEF Guía de Uso"
# Mark the "GUIAS DE USO" status column (H) for the Coeficientes PEF rows as
# "ACTUALIZADO", and mark a handful of other rows (module guides that don't
# apply) as "N/A".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coeficientes block: PEF / AJUSTE ANUAL / AJUSTE SEMESTRAL -> ACTUALIZADO
$ws.Range("H36").Value = "ACTUALIZADO"
$ws.Range("H37").Value = "ACTUALIZADO"
$ws.Range("H38").Value = "ACTUALIZADO"

# Módulo DAMOP rows (CARGA DE SPEI Y CFDI ... SOLICITUDES DE ANTICIPOS) -> N/A
$ws.Range("H66").Value = "N/A"
$ws.Range("H67").Value = "N/A"
$ws.Range("H68").Value = "N/A"
$ws.Range("H69").Value = "N/A"
$ws.Range("H70").Value = "N/A"

# CONSULTA DE SPEI / CARGA DE CFDI rows -> N/A, matching the shading of the
# rest of that merged block (copy format from a sibling row, then set value).
[void]$ws.Range("H79").Copy()
$ws.Range("H85").PasteSpecial(-4122)
$ws.Range("H86").PasteSpecial(-4122)
$ws.Range("H85").Value = "N/A"
$ws.Range("H86").Value = "N/A"

# Leave the cursor where the author last left it when saving.
[void]$ws.Range("J79").Select()
